$wb = $excel.ActiveWorkbook

# ---- Caso1 ----
$ws = $wb.Worksheets.Item("Caso1")

# Header row additions
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"

# Data rows: update column G, add columns H, I, J
$ws.Range("G2").Value = 0.9975649118423462
$ws.Range("H2").Value = 0.9937247037887573
$ws.Range("I2").Value = 0.9942487114884958
$ws.Range("J2").Value = 0.9939124397933483
$ws.Range("G3").Value = 0.9940235018730164
$ws.Range("H3").Value = 0.9937007427215576
$ws.Range("I3").Value = 0.9937017374150956
$ws.Range("J3").Value = 0.9935125056654215
$ws.Range("G4").Value = 0.994796097278595
$ws.Range("H4").Value = 0.9937124252319336
$ws.Range("I4").Value = 0.9938039122923608
$ws.Range("J4").Value = 0.9936844762414694
$ws.Range("G5").Value = 0.993896484375
$ws.Range("H5").Value = 0.993661642074585
$ws.Range("I5").Value = 0.9937972011340818
$ws.Range("J5").Value = 0.9937595035880804
$ws.Range("G6").Value = 0.9941811561584473
$ws.Range("H6").Value = 0.9934757351875305
$ws.Range("I6").Value = 0.9938649233596956
$ws.Range("J6").Value = 0.9938138220459223
$ws.Range("G7").Value = 0.9938853979110718
$ws.Range("H7").Value = 0.9937895536422729
$ws.Range("I7").Value = 0.993877168238681
$ws.Range("J7").Value = 0.9937603436410427
$ws.Range("G8").Value = 0.9940674901008606
$ws.Range("H8").Value = 0.9936531186103821
$ws.Range("I8").Value = 0.9938863441752228
$ws.Range("J8").Value = 0.9937184043228626
$ws.Range("G9").Value = 0.9943283200263977
$ws.Range("H9").Value = 0.9934903979301453
$ws.Range("I9").Value = 0.9938985135539904
$ws.Range("J9").Value = 0.9938072059303522
$ws.Range("G10").Value = 0.9946186542510986
$ws.Range("H10").Value = 0.9934080839157104
$ws.Range("I10").Value = 0.9939343274396578
$ws.Range("J10").Value = 0.9937522038817406
$ws.Range("G11").Value = 0.9938380718231201
$ws.Range("H11").Value = 0.9935124516487122
$ws.Range("I11").Value = 0.9939341158998262
$ws.Range("J11").Value = 0.9938011411577463
$ws.Range("G12").Value = 0.9946953058242798
$ws.Range("H12").Value = 0.9936211705207825
$ws.Range("I12").Value = 0.9939699313552328
$ws.Range("J12").Value = 0.9938318245112896
$ws.Range("G13").Value = 0.9938383102416992
$ws.Range("H13").Value = 0.9937409162521362
$ws.Range("I13").Value = 0.9939721863074784
$ws.Range("J13").Value = 0.9938335344195366
$ws.Range("G14").Value = 0.9975169897079468
$ws.Range("H14").Value = 0.9950545430183411
$ws.Range("I14").Value = 0.9949798073298636
$ws.Range("J14").Value = 0.99478062056005
$ws.Range("G15").Value = 0.9985942840576172
$ws.Range("H15").Value = 0.9947940111160278
$ws.Range("I15").Value = 0.9950405652747434
$ws.Range("J15").Value = 0.9948516208678484
$ws.Range("G16").Value = 0.9941329956054688
$ws.Range("H16").Value = 0.9936136603355408
$ws.Range("I16").Value = 0.9937832682597182
$ws.Range("J16").Value = 0.9936010111123323
$ws.Range("G17").Value = 0.9944035410881042
$ws.Range("H17").Value = 0.9934131503105164
$ws.Range("I17").Value = 0.9937726762869395
$ws.Range("J17").Value = 0.9935921411961317
$ws.Range("G18").Value = 0.9944289326667786
$ws.Range("H18").Value = 0.9934409856796265
$ws.Range("I18").Value = 0.9937853315710676
$ws.Range("J18").Value = 0.9936011023819447
$ws.Range("G19").Value = 0.9937174916267395
$ws.Range("H19").Value = 0.9933353662490845
$ws.Range("I19").Value = 0.993779603544176
$ws.Range("J19").Value = 0.9935773424804211

# ---- Caso2 ----
$ws = $wb.Worksheets.Item("Caso2")

# Header row additions
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"

# Data rows: update column G, add columns H, I, J
$ws.Range("G2").Value = 0.9963275194168091
$ws.Range("H2").Value = 0.9922446012496948
$ws.Range("I2").Value = 0.9927920481799234
$ws.Range("J2").Value = 0.9923256225883961
$ws.Range("G3").Value = 0.9926695227622986
$ws.Range("H3").Value = 0.9919544458389282
$ws.Range("I3").Value = 0.99222456725438
$ws.Range("J3").Value = 0.9919021427631378
$ws.Range("G4").Value = 0.9933093786239624
$ws.Range("H4").Value = 0.9917682409286499
$ws.Range("I4").Value = 0.9923333267797528
$ws.Range("J4").Value = 0.9920868910849094
$ws.Range("G5").Value = 0.9924783706665039
$ws.Range("H5").Value = 0.9917528629302979
$ws.Range("I5").Value = 0.9923269832824668
$ws.Range("J5").Value = 0.9921437371522188
$ws.Range("G6").Value = 0.9927784204483032
$ws.Range("H6").Value = 0.991658091545105
$ws.Range("I6").Value = 0.9923970721791174
$ws.Range("J6").Value = 0.9922917298972607
$ws.Range("G7").Value = 0.9924938678741455
$ws.Range("H7").Value = 0.9919619560241699
$ws.Range("I7").Value = 0.9924088460429616
$ws.Range("J7").Value = 0.9922812953591347
$ws.Range("G8").Value = 0.9926480054855347
$ws.Range("H8").Value = 0.9918528199195862
$ws.Range("I8").Value = 0.992419658051451
$ws.Range("J8").Value = 0.9922273512929678
$ws.Range("G9").Value = 0.9929164052009583
$ws.Range("H9").Value = 0.9917550086975098
$ws.Range("I9").Value = 0.9924319278717652
$ws.Range("J9").Value = 0.9923817403614521
$ws.Range("G10").Value = 0.9931668043136597
$ws.Range("H10").Value = 0.991756796836853
$ws.Range("I10").Value = 0.9924656892219672
$ws.Range("J10").Value = 0.9923104345798492
$ws.Range("G11").Value = 0.9924137592315674
$ws.Range("H11").Value = 0.9917746782302856
$ws.Range("I11").Value = 0.9924672113235948
$ws.Range("J11").Value = 0.992219515144825
$ws.Range("G12").Value = 0.9932748675346375
$ws.Range("H12").Value = 0.9919631481170654
$ws.Range("I12").Value = 0.9925047477988744
$ws.Range("J12").Value = 0.992371516302228
$ws.Range("G13").Value = 0.9924613237380981
$ws.Range("H13").Value = 0.9920219779014587
$ws.Range("I13").Value = 0.9925087652560451
$ws.Range("J13").Value = 0.992353668436408
$ws.Range("G14").Value = 0.996087908744812
$ws.Range("H14").Value = 0.9932867884635925
$ws.Range("I14").Value = 0.9935337835709764
$ws.Range("J14").Value = 0.9933609012514353
$ws.Range("G15").Value = 0.9971036314964294
$ws.Range("H15").Value = 0.9930979609489441
$ws.Range("I15").Value = 0.9935941134127776
$ws.Range("J15").Value = 0.9934268519282341
$ws.Range("G16").Value = 0.9927046895027161
$ws.Range("H16").Value = 0.991805374622345
$ws.Range("I16").Value = 0.9923124706973084
$ws.Range("J16").Value = 0.9920274298638105
$ws.Range("G17").Value = 0.9929711818695068
$ws.Range("H17").Value = 0.9916262030601501
$ws.Range("I17").Value = 0.9923024928311004
$ws.Range("J17").Value = 0.9920345675200224
$ws.Range("G18").Value = 0.9929564595222473
$ws.Range("H18").Value = 0.9915755987167358
$ws.Range("I18").Value = 0.9923152392004456
$ws.Range("J18").Value = 0.9919967800378799
$ws.Range("G19").Value = 0.992320716381073
$ws.Range("H19").Value = 0.9914385080337524
$ws.Range("I19").Value = 0.992306410030156
$ws.Range("J19").Value = 0.9920419789850712

# ---- Caso3 ----
$ws = $wb.Worksheets.Item("Caso3")

# Header row additions
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"

# Data rows: update column G, add columns H, I, J
$ws.Range("G2").Value = 0.9947376847267151
$ws.Range("H2").Value = 0.990793764591217
$ws.Range("I2").Value = 0.9910289857145024
$ws.Range("J2").Value = 0.9908338226377964
$ws.Range("G3").Value = 0.9911172389984131
$ws.Range("H3").Value = 0.9907561540603638
$ws.Range("I3").Value = 0.9905169528441964
$ws.Range("J3").Value = 0.9904487542808056
$ws.Range("G4").Value = 0.991762638092041
$ws.Range("H4").Value = 0.990637481212616
$ws.Range("I4").Value = 0.9906388149055246
$ws.Range("J4").Value = 0.9906573742628098
$ws.Range("G5").Value = 0.9909349083900452
$ws.Range("H5").Value = 0.990558922290802
$ws.Range("I5").Value = 0.9906353328603542
$ws.Range("J5").Value = 0.9906699694693089
$ws.Range("G6").Value = 0.9912029504776001
$ws.Range("H6").Value = 0.9905632734298706
$ws.Range("I6").Value = 0.9907105675812984
$ws.Range("J6").Value = 0.9907280560582876
$ws.Range("G7").Value = 0.990939199924469
$ws.Range("H7").Value = 0.9908116459846497
$ws.Range("I7").Value = 0.99072232141887
$ws.Range("J7").Value = 0.9907101131975651
$ws.Range("G8").Value = 0.9911137223243713
$ws.Range("H8").Value = 0.99065101146698
$ws.Range("I8").Value = 0.9907329275712168
$ws.Range("J8").Value = 0.9907309729605913
$ws.Range("G9").Value = 0.991344153881073
$ws.Range("H9").Value = 0.990490198135376
$ws.Range("I9").Value = 0.9907433571308604
$ws.Range("J9").Value = 0.9907716736197472
$ws.Range("G10").Value = 0.9916150569915771
$ws.Range("H10").Value = 0.9904554486274719
$ws.Range("I10").Value = 0.9907832359875928
$ws.Range("J10").Value = 0.9907750263810158
$ws.Range("G11").Value = 0.9908909797668457
$ws.Range("H11").Value = 0.9904032945632935
$ws.Range("I11").Value = 0.9907840527966344
$ws.Range("J11").Value = 0.9907806944102049
$ws.Range("G12").Value = 0.9916942119598389
$ws.Range("H12").Value = 0.9906036257743835
$ws.Range("I12").Value = 0.9908205630505528
$ws.Range("J12").Value = 0.9908340889960527
$ws.Range("G13").Value = 0.9909060001373291
$ws.Range("H13").Value = 0.9906781911849976
$ws.Range("I13").Value = 0.99082368362839
$ws.Range("J13").Value = 0.99087786488235
$ws.Range("G14").Value = 0.994605541229248
$ws.Range("H14").Value = 0.99207603931427
$ws.Range("I14").Value = 0.9918727743215692
$ws.Range("J14").Value = 0.991892758756876
$ws.Range("G15").Value = 0.9956146478652954
$ws.Range("H15").Value = 0.9918802380561829
$ws.Range("I15").Value = 0.9919345928434188
$ws.Range("J15").Value = 0.991925872862339
$ws.Range("G16").Value = 0.9911488890647888
$ws.Range("H16").Value = 0.9905741810798645
$ws.Range("I16").Value = 0.9906203502513335
$ws.Range("J16").Value = 0.990576907992363
$ws.Range("G17").Value = 0.9913815855979919
$ws.Range("H17").Value = 0.9904745817184448
$ws.Range("I17").Value = 0.9906112752181409
$ws.Range("J17").Value = 0.9905827473849058
$ws.Range("G18").Value = 0.9914137721061707
$ws.Range("H18").Value = 0.9903580546379089
$ws.Range("I18").Value = 0.9906207553786884
$ws.Range("J18").Value = 0.99054323323071
$ws.Range("G19").Value = 0.9907498955726624
$ws.Range("H19").Value = 0.9903942346572876
$ws.Range("I19").Value = 0.9906154789465073
$ws.Range("J19").Value = 0.990569407120347

Write-Output "Edit applied successfully"
